$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 25.91000000000061
$ws.Range("G2").Value = 0.002106471215033379
$ws.Range("H2").Value = 0.007508572287002756
$ws.Range("K2").Value = 4.200872594137925
$ws.Range("L2").Value = "[1.0104481374449747, 7.391297050830875]"
$ws.Range("M2").Value = 0.010015678487997
$ws.Range("N2").Value = 0.010015678487997
$ws.Range("O2").Value = -0.9182633181663862
$ws.Range("P2").Value = "[-1.610105544182157, -0.22642109215061534]"
$ws.Range("Q2").Value = 0.009438904889965904
$ws.Range("R2").Value = 0.009438904889965904
$ws.Range("S2").Value = 13.16269639131862
$ws.Range("T2").Value = "[11.51245743705022, 14.812935345587015]"
$ws.Range("W2").Value = 3.786646646646737
$ws.Range("X2").Value = 0.9336936936937144
$ws.Range("Y2").Value = 6.639599599599761

# Row 3
$ws.Range("E3").Value = 24.03000000000032
$ws.Range("G3").Value = [double]"9.653556762145854e-05"
$ws.Range("H3").Value = 0.001631817327551451
$ws.Range("K3").Value = 4.768993590219882
$ws.Range("L3").Value = "[2.1355497806312442, 7.40243739980852]"
$ws.Range("M3").Value = 0.0004132301852455278
$ws.Range("N3").Value = 0.0008264603704910556
$ws.Range("O3").Value = 2.899447874484274
$ws.Range("P3").Value = "[2.245342497160273, 3.5535532518082746]"
$ws.Range("S3").Value = 13.85097417045925
$ws.Range("T3").Value = "[12.327696567556677, 15.374251773361815]"
$ws.Range("W3").Value = 12.94108108108125
$ws.Range("X3").Value = 10.4394594594596
$ws.Range("Y3").Value = 15.44270270270291
